$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.259.72'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.91%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.719.96'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.39%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.88'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.85%  '

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4708'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.00%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2625'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06204'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.73%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.717.25'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07071'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.42%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.33'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.81%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6013'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.47%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.431'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.40%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.20'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.80%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.02%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.02%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.271.87'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.97%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006804'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.05%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.56'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.76%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.935.57'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.27%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.39%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.728'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.51%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.289'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.19%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.72'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.48%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.17'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.96%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.47%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.97%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '107.16'

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.974'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.27%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.44%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.07770'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.22%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04459'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.93%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.616'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.22%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9762'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.87%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6181'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9380'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +8.24%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '112.54'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +16.38%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.440'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.29%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.926'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +3.84%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.0000'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.06%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01479'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.445'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +12.94%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3824'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.62%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1179'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.79%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.280'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.23%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05268'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.786'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +6.48%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.23'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.87%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3381'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.44%  '

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.216'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.53%  '
